$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as exact text (avoid float auto-conversion
# of numeric-looking strings such as "233.20" -> 233.2).
$ws.Range("D2").Value = '37.808.69'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.078.35'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.20'
$ws.Range("E5").Value = '  -0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.626'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.84'
$ws.Range("E7").Value = '  -1.32%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.393'
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0787'
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("E11").Value = '  +3.64%  '
$ws.Range("D12").Value = '2.383.94'
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.80'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.15'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.785'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.36'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '2.048.68'
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").Value = '37.694.79'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").Value = '  -1.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.63'
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '0.0₃0843'
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.71'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.77'
$ws.Range("E26").Value = '  +7.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.77'
$ws.Range("E27").Value = '  +0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.139'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("E29").Value = '  -2.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.44'
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.122'
$ws.Range("E31").Value = '  +1.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.75'
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.71'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.83'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.42'
$ws.Range("E37").Value = '  -3.22%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("E40").Value = '  +8.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.67'
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0978'
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.20'
$ws.Range("E43").Value = '  +6.44%  '
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '1.450.63'
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.07'
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.10'
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.40'
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("D51").Value = '2.269.02'
$ws.Range("E51").Value = '  -0.43%  '
